{"js": "// Change \"could fade\" -> \"could die down\", while preserving the original\n// two-run split (run 1 keeps its own formatting and becomes \"could \",\n// run 2 keeps its own formatting and becomes \"die down\").\n\nconst body = context.document.body;\n\n// Locate the phrase as a whole (it is rendered contiguously even though it\n// is stored as two separate runs: \"could fa\" + \"de\").\nconst matches = body.search(\"could fade\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find \"could fade\" in the document body.');\n}\n\nconst whole = matches.items[0];\n\n// Re-search, but scoped to just this occurrence, so we can independently\n// target each of the two runs that make up \"could fade\" without touching\n// any other occurrence of the short substring \"de\" elsewhere in the doc.\nconst firstRunMatches = whole.search(\"could fa\", { matchCase: true });\nfirstRunMatches.load(\"items\");\nawait context.sync();\n\nif (firstRunMatches.items.length !== 1) {\n  throw new Error('Expected exactly one \"could fa\" match within the located phrase.');\n}\n\n// Replace text in-place; since this range corresponds exactly to the first\n// run, Word keeps that run's own run properties (rPr) intact.\nfirstRunMatches.items[0].insertText(\"could \", \"Replace\");\n\nconst secondRunMatches = whole.search(\"de\", { matchCase: true });\nsecondRunMatches.load(\"items\");\nawait context.sync();\n\nif (secondRunMatches.items.length !== 1) {\n  throw new Error('Expected exactly one \"de\" match within the located phrase.');\n}\n\n// Replace text in-place; this range corresponds exactly to the second run,\n// so its own (different) run properties are preserved as well.\nsecondRunMatches.items[0].insertText(\"die down\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Change \"could fade\" -> \"could die down\", while preserving the original\n# two-run split (run 1 keeps its own formatting and becomes \"could \",\n# run 2 keeps its own formatting and becomes \"die down\").\n\n$d = $word.ActiveDocument\n\n$oldFirst  = \"could fa\"\n$oldSecond = \"de\"\n$newFirst  = \"could \"\n$newSecond = \"die down\"\n$wholeOld  = $oldFirst + $oldSecond   # \"could fade\" (rendered contiguously,\n                                      # even though stored as two runs)\n\n# Sanity check: make sure the phrase is unique in the document so we edit\n# the right (and only the right) spot.\n$fullText = $d.Content.Text\n$count = 0\n$searchIdx = -1\nwhile ($true) {\n  $searchIdx = $fullText.IndexOf($wholeOld, $searchIdx + 1)\n  if ($searchIdx -eq -1) { break }\n  $count++\n}\nif ($count -ne 1) {\n  throw \"expected exactly one occurrence of '$wholeOld', found $count\"\n}\n\n# Step 1: locate the phrase precisely using Find, scoped to the whole\n# document content.\n$whole = $d.Content\n$foundWhole = $whole.Find.Execute($wholeOld, $true, $false, $false, $false, $false, $true, 1)\nif (-not $foundWhole) {\n  throw \"could not find '$wholeOld'\"\n}\n$wholeStart = $whole.Start\n\n# Step 2: replace the first run's text (\"could fa\" -> \"could \"). Scope the\n# Find tightly to the located phrase, then assign .Text directly so the\n# run keeps its own run properties (rPr) untouched.\n$r1 = $d.Range($wholeStart, $wholeStart + $wholeOld.Length)\n$found1 = $r1.Find.Execute($oldFirst, $true)\nif (-not $found1) {\n  throw \"could not find '$oldFirst' within the located phrase\"\n}\n$r1.Text = $newFirst\n\n# Step 3: replace the second run's text (\"de\" -> \"die down\"). Re-scope the\n# range tightly to right after the newly written first-run text, so this\n# cannot accidentally match any other \"de\" substring elsewhere in the\n# document (e.g. inside \"side\", \"del\", etc.).\n$r2Start = $wholeStart + $newFirst.Length\n$r2 = $d.Range($r2Start, $r2Start + $oldSecond.Length)\n$found2 = $r2.Find.Execute($oldSecond, $true)\nif (-not $found2) {\n  throw \"could not find '$oldSecond' within the located phrase\"\n}\n$r2.Text = $newSecond\n"}
